# Update countries & provincias Spain
# Refresh COVID-19 country stats + re-sort consequences (a few countries
# swap adjacent rows because the table is kept sorted by "Casos totales"
# descending), and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp -------------------------------------
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 22 de Mayo de 2020 a las 20:05"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Cells.Item(4, 2).Value2 = 1632629
$ws.Cells.Item(4, 3).Value2 = 11727
$ws.Cells.Item(4, 5).Value2 = 1150457
$ws.Cells.Item(4, 7).Value2 = 748
$ws.Cells.Item(4, 8).Value2 = 97102

# --- India (row 14) ----------------------------------------------------------
$ws.Cells.Item(14, 2).Value2 = 124462
$ws.Cells.Item(14, 3).Value2 = 6236
$ws.Cells.Item(14, 4).Value2 = 51687
$ws.Cells.Item(14, 5).Value2 = 69063
$ws.Cells.Item(14, 7).Value2 = 128
$ws.Cells.Item(14, 8).Value2 = 3712

# --- Sudafrica (row 37) -------------------------------------------------------
$ws.Cells.Item(37, 4).Value2 = 10104
$ws.Cells.Item(37, 5).Value2 = 9624
$ws.Cells.Item(37, 7).Value2 = 28
$ws.Cells.Item(37, 8).Value2 = 397

# --- Oman / Moldavia swap (rows 62-63) ---------------------------------------
$ws.Cells.Item(62, 1).Value2 = "Moldavia"
$ws.Cells.Item(62, 2).Value2 = 6847
$ws.Cells.Item(62, 3).Value2 = 143
$ws.Cells.Item(62, 4).Value2 = 3369
$ws.Cells.Item(62, 5).Value2 = 3241
$ws.Cells.Item(62, 7).Value2 = 4
$ws.Cells.Item(62, 8).Value2 = 237

$ws.Cells.Item(63, 1).Value2 = "Oman"
$ws.Cells.Item(63, 2).Value2 = 6794
$ws.Cells.Item(63, 3).Value2 = 424
$ws.Cells.Item(63, 4).Value2 = 1821
$ws.Cells.Item(63, 5).Value2 = 4941
$ws.Cells.Item(63, 7).Value2 = 1
$ws.Cells.Item(63, 8).Value2 = 32

# --- Republica del Chad moves ahead of Sierra Leona / Malta (rows 126-128) ---
$ws.Cells.Item(126, 1).Value2 = "Republica del Chad"
$ws.Cells.Item(126, 2).Value2 = 611
$ws.Cells.Item(126, 3).Value2 = 23
$ws.Cells.Item(126, 4).Value2 = 196
$ws.Cells.Item(126, 5).Value2 = 357
$ws.Cells.Item(126, 7).Value2 = 0
$ws.Cells.Item(126, 8).Value2 = 58

$ws.Cells.Item(127, 1).Value2 = "Sierra Leona"
$ws.Cells.Item(127, 2).Value2 = 606
$ws.Cells.Item(127, 3).Value2 = 21
$ws.Cells.Item(127, 4).Value2 = 230
$ws.Cells.Item(127, 5).Value2 = 338
$ws.Cells.Item(127, 7).Value2 = 3
$ws.Cells.Item(127, 8).Value2 = 38

$ws.Cells.Item(128, 1).Value2 = "Malta"
$ws.Cells.Item(128, 2).Value2 = 600
$ws.Cells.Item(128, 3).Value2 = 1
$ws.Cells.Item(128, 4).Value2 = 469
$ws.Cells.Item(128, 5).Value2 = 125
$ws.Cells.Item(128, 8).Value2 = 6

# --- Ruanda (row 146) ---------------------------------------------------------
$ws.Cells.Item(146, 2).Value2 = 321
$ws.Cells.Item(146, 3).Value2 = 1
$ws.Cells.Item(146, 4).Value2 = 222
$ws.Cells.Item(146, 5).Value2 = 99

# --- Yemen / Birmania swap (rows 152-153) -------------------------------------
$ws.Cells.Item(152, 1).Value2 = "Yemen"
$ws.Cells.Item(152, 2).Value2 = 209
$ws.Cells.Item(152, 3).Value2 = 12
$ws.Cells.Item(152, 4).Value2 = 11
$ws.Cells.Item(152, 5).Value2 = 165
$ws.Cells.Item(152, 8).Value2 = 33

$ws.Cells.Item(153, 1).Value2 = "Birmania"
$ws.Cells.Item(153, 2).Value2 = 199
$ws.Cells.Item(153, 4).Value2 = 108
$ws.Cells.Item(153, 5).Value2 = 85
$ws.Cells.Item(153, 8).Value2 = 6

# --- Lesoto (row 218) ----------------------------------------------------------
$ws.Cells.Item(218, 2).Value2 = 2
$ws.Cells.Item(218, 3).Value2 = 1
$ws.Cells.Item(218, 5).Value2 = 2
